# Updated cryptos list values per target diff (rows 2-51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''66.848.73'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.49%  '
$ws.Range("D3").Value = '''3.497.34'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.32%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '''595.07'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.75%  '
$ws.Range("D6").Value = '''172.69'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.49%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").Value = '''0.132'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.13%  '
$ws.Range("D10").Value = '''7.17'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.94%  '
$ws.Range("E11").Value = '  -0.95%  '
$ws.Range("D12").Value = '''4.099.13'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.24%  '
$ws.Range("E13").Value = '  -0.09%  '
$ws.Range("D14").Value = '''29.34'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.32%  '
$ws.Range("D15").Value = '''66.861.38'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.47%  '
$ws.Range("E16").Value = '  +0.69%  '
$ws.Range("D17").Value = '''3.490.10'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.19%  '
$ws.Range("D18").Value = '''6.28'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.24%  '
$ws.Range("E19").Value = '  +1.87%  '
$ws.Range("D20").Value = '''394.67'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.08%  '
$ws.Range("D21").Value = '''7.92'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.23%  '
$ws.Range("D22").Value = '''73.25'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.46%  '
$ws.Range("D23").Value = '''1.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.13%  '
$ws.Range("D24").Value = '''0.534'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.39%  '
$ws.Range("E25").Value = '  +0.15%  '
$ws.Range("D26").Value = '''10.20'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.50%  '
$ws.Range("E27").Value = '  +0.51%  '
$ws.Range("D28").Value = '''0.988'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.31%  '
$ws.Range("E29").Value = '  -2.57%  '
$ws.Range("E30").Value = '  -1.57%  '
$ws.Range("E31").Value = '  +0.03%  '
$ws.Range("D32").Value = '''23.70'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.68%  '
$ws.Range("E33").Value = '  -0.20%  '
$ws.Range("E34").Value = '  +0.44%  '
$ws.Range("D35").Value = '''162.43'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.11%  '
$ws.Range("D36").Value = '''0.876'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.25%  '
$ws.Range("D37").Value = '''1.90'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.61%  '
$ws.Range("D38").Value = '''6.88'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.57%  '
$ws.Range("E39").Value = '  +0.26%  '
$ws.Range("E40").Value = '  -0.73%  '
$ws.Range("B41").Value = 'Maker'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D41").Value = '''2.833.71'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.50%  '
$ws.Range("B42").Value = 'InjectiveProtocol'
$ws.Range("C42").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D42").Value = '''27.13'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.62%  '
$ws.Range("D43").Value = '''26.12'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.60%  '
$ws.Range("E44").Value = '  -0.88%  '
$ws.Range("E45").Value = '  +2.06%  '
$ws.Range("E46").Value = '  -2.90%  '
$ws.Range("D47").Value = '''337.21'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.17%  '
$ws.Range("D48").Value = '''34.59'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.23%  '
$ws.Range("E49").Value = '  -1.11%  '
$ws.Range("E50").Value = '  -1.10%  '
$ws.Range("D51").Value = '''0.841'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.47%  '
